$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.792.55'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '1.619.21'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.79%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.991'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.83%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.21'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.27%  '
$ws.Range('E9').Value = '  +3.16%  '
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '1.851.61'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = '1.631.34'
$ws.Range('E13').Value = '  +1.39%  '
$ws.Range('E14').Value = '  +5.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.89'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.91%  '
$ws.Range('D16').Value = '29.828.53'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +16.96%  '
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.25'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.992'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('E22').Value = '  +2.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.64'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.64%  '
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +3.19%  '
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('E30').Value = '  +2.85%  '
$ws.Range('E31').Value = '  +2.96%  '
$ws.Range('E32').Value = '  +3.38%  '
$ws.Range('E33').Value = '  +3.76%  '
$ws.Range('D34').Value = '1.419.72'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('E35').Value = '  +6.53%  '
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.28'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('E39').Value = '  +2.36%  '
$ws.Range('E40').Value = '  +3.86%  '
$ws.Range('E41').Value = '  +3.47%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.97'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.24%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.826'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '53.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.02'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +18.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.990'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('E48').Value = '  +2.96%  '
$ws.Range('D49').Value = '1.760.51'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '88.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0533'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.98%  '
